# Update a handful of numeric values in Sheet1 to reflect the revised
# RandomForest imputation results ("Update Name of Algo" commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value  = -21.43140000000002
$ws.Range("D5").Value  = -8.379299999999995
$ws.Range("E7").Value  = 12.0994
$ws.Range("D9").Value  = -8.521400000000002
$ws.Range("D11").Value = -8.325700000000008
$ws.Range("E11").Value = 13.08889999999999
$ws.Range("A21").Value = -21.35930000000001
$ws.Range("D21").Value = -7.880200000000002
$ws.Range("E21").Value = 12.98680000000001
$ws.Range("A23").Value = -21.40460000000003
$ws.Range("A25").Value = -22.45000000000003
